# Updated cryptos list (price/volume refresh) - mirrors the GitHub Actions data pull.
# Cells whose new text could otherwise be auto-parsed as a number by Excel
# are written with a leading apostrophe and then restored to the default
# "Normal" style so the resulting cell stays plain text with no style index,
# matching the original workbook's formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.766.08"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "1.594.63"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("D5").Value = "'210.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'19.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("D12").Value = "1.818.47"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").Value = "1.593.69"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "'0.528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "26.741.08"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "'63.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "'209.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.43%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("E23").Value = "  -6.88%  "
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").Value = "'146.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'0.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("D29").Value = "'15.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "'0.0500"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").Value = "'0.671"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +23.93%  "
$ws.Range("D34").Value = "'2.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").Value = "1.310.03"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  -3.38%  "
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "'0.818"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "'62.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("D45").Value = "1.730.89"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").Value = "'89.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").Value = "'1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'0.807"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -4.41%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0509"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0973"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.42%  "
